$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4150643333333333
$ws.Range("H2").Value = 1.245193
$ws.Range("I2").Value = 0.02396302145531912
$ws.Range("J2").Value = 0.02396302145531911
$ws.Range("M2").Value = 1.0848515
$ws.Range("N2").Value = 2.169703
$ws.Range("O2").Value = 0.2186227527895346
$ws.Range("P2").Value = 0.1768985792936573
$ws.Range("Q2").Value = 0.4502831646131667
$ws.Range("R2").Value = 2.701698987679
$ws.Range("S2").Value = 0.005238861715716545
$ws.Range("T2").Value = 0.004239024451029381

# Row 3
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4150643333333333
$ws.Range("H3").Value = 1.245193
$ws.Range("I3").Value = 0.02396302145531912
$ws.Range("J3").Value = 0.02396302145531911
$ws.Range("O3").Value = 0.2203773817607929
$ws.Range("P3").Value = 0.267477505734296
$ws.Range("Q3").Value = 0.4538970605861111
$ws.Range("R3").Value = 4.085073545275
$ws.Range("S3").Value = 0.005280907927400932
$ws.Range("T3").Value = 0.006409569208726176

# Row 4
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.4150643333333333
$ws.Range("H4").Value = 1.245193
$ws.Range("I4").Value = 0.02396302145531912
$ws.Range("J4").Value = 0.02396302145531911
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1406063333333333
$ws.Range("N4").Value = 0.421819
$ws.Range("O4").Value = 0.02833543913888328
$ws.Range("P4").Value = 0.03439142676166795
$ws.Range("Q4").Value = 0.05836067400744444
$ws.Range("R4").Value = 0.525246066067
$ws.Range("S4").Value = 0.000679002736030949
$ws.Range("T4").Value = 0.0008241224973688851

# Row 5
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.4150643333333333
$ws.Range("H5").Value = 1.245193
$ws.Range("I5").Value = 0.02396302145531912
$ws.Range("J5").Value = 0.02396302145531911
$ws.Range("M5").Value = 1.5365345
$ws.Range("N5").Value = 3.073069
$ws.Range("O5").Value = 0.3096473592432615
$ws.Range("P5").Value = 0.2505511308097838
$ws.Range("Q5").Value = 0.6377606678861667
$ws.Range("R5").Value = 3.826564007317
$ws.Range("S5").Value = 0.007420086313129183
$ws.Range("T5").Value = 0.006003962123249315

# Row 6
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.4150643333333333
$ws.Range("H6").Value = 1.245193
$ws.Range("I6").Value = 0.02396302145531912
$ws.Range("J6").Value = 0.02396302145531911
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3224876666666667
$ws.Range("N6").Value = 0.967463
$ws.Range("O6").Value = 0.06498874862351253
$ws.Range("P6").Value = 0.07887845950306545
$ws.Range("Q6").Value = 0.1338531283732222
$ws.Range("R6").Value = 1.204678155359
$ws.Range("S6").Value = 0.001557326777619572
$ws.Range("T6").Value = 0.001890166217434477

# Row 7
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.4150643333333333
$ws.Range("H7").Value = 1.245193
$ws.Range("I7").Value = 0.02396302145531912
$ws.Range("J7").Value = 0.02396302145531911
$ws.Range("M7").Value = 0.7841693333333334
$ws.Range("N7").Value = 2.352508
$ws.Range("O7").Value = 0.1580283184440151
$ws.Range("P7").Value = 0.1918028978975294
$ws.Range("Q7").Value = 0.3254807215604444
$ws.Range("R7").Value = 2.929326494044
$ws.Range("S7").Value = 0.003786835985421936
$ws.Range("T7").Value = 0.004596176957510878

# Row 8
$ws.Range("A8").Value = "Neutrophils"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 16.787094
$ws.Range("H8").Value = 50.361282
$ws.Range("I8").Value = 0.9691738397849783
$ws.Range("J8").Value = 0.9691738397849782
$ws.Range("M8").Value = 1.0848515
$ws.Range("N8").Value = 2.169703
$ws.Range("O8").Value = 0.2186227527895346
$ws.Range("P8").Value = 0.1768985792936573
$ws.Range("Q8").Value = 18.211504106541
$ws.Range("R8").Value = 109.269024639246
$ws.Range("S8").Value = 0.2118834527853953
$ws.Range("T8").Value = 0.1714454753465413

# Row 9
$ws.Range("A9").Value = "Neutrophils"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 16.787094
$ws.Range("H9").Value = 50.361282
$ws.Range("I9").Value = 0.9691738397849783
$ws.Range("J9").Value = 0.9691738397849782
$ws.Range("O9").Value = 0.2203773817607929
$ws.Range("P9").Value = 0.267477505734296
$ws.Range("Q9").Value = 18.35766653615
$ws.Range("R9").Value = 165.21899882535
$ws.Range("S9").Value = 0.2135839932828677
$ws.Range("T9").Value = 0.2592322012886162

# Row 10
$ws.Range("A10").Value = "Neutrophils"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 16.787094
$ws.Range("H10").Value = 50.361282
$ws.Range("I10").Value = 0.9691738397849783
$ws.Range("J10").Value = 0.9691738397849782
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.1406063333333333
$ws.Range("N10").Value = 0.421819
$ws.Range("O10").Value = 0.02833543913888328
$ws.Range("P10").Value = 0.03439142676166795
$ws.Range("Q10").Value = 2.360371734662
$ws.Range("R10").Value = 21.243345611958
$ws.Range("S10").Value = 0.02746196635222506
$ws.Range("T10").Value = 0.03333127113028959

# Row 11
$ws.Range("A11").Value = "Neutrophils"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 16.787094
$ws.Range("H11").Value = 50.361282
$ws.Range("I11").Value = 0.9691738397849783
$ws.Range("J11").Value = 0.9691738397849782
$ws.Range("M11").Value = 1.5365345
$ws.Range("N11").Value = 3.073069
$ws.Range("O11").Value = 0.3096473592432615
$ws.Range("P11").Value = 0.2505511308097838
$ws.Range("Q11").Value = 25.793949085743
$ws.Range("R11").Value = 154.763694514458
$ws.Range("S11").Value = 0.3001021201370704
$ws.Range("T11").Value = 0.2428276015093865

# Row 12
$ws.Range("A12").Value = "Neutrophils"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 16.787094
$ws.Range("H12").Value = 50.361282
$ws.Range("I12").Value = 0.9691738397849783
$ws.Range("J12").Value = 0.9691738397849782
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.3224876666666667
$ws.Range("N12").Value = 0.967463
$ws.Range("O12").Value = 0.06498874862351253
$ws.Range("P12").Value = 0.07887845950306545
$ws.Range("Q12").Value = 5.413630774174
$ws.Range("R12").Value = 48.722676967566
$ws.Range("S12").Value = 0.06298539504627036
$ws.Range("T12").Value = 0.07644693947290984

# Row 13
$ws.Range("A13").Value = "Neutrophils"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 16.787094
$ws.Range("H13").Value = 50.361282
$ws.Range("I13").Value = 0.9691738397849783
$ws.Range("J13").Value = 0.9691738397849782
$ws.Range("M13").Value = 0.7841693333333334
$ws.Range("N13").Value = 2.352508
$ws.Range("O13").Value = 0.1580283184440151
$ws.Range("P13").Value = 0.1918028978975294
$ws.Range("Q13").Value = 13.163924310584
$ws.Range("R13").Value = 118.475318795256
$ws.Range("S13").Value = 0.1531569121811495
$ws.Range("T13").Value = 0.1858903510372347

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1188766666666667
$ws.Range("H14").Value = 0.35663
$ws.Range("I14").Value = 0.006863138759702679
$ws.Range("J14").Value = 0.006863138759702677
$ws.Range("M14").Value = 1.0848515
$ws.Range("N14").Value = 2.169703
$ws.Range("O14").Value = 0.2186227527895346
$ws.Range("P14").Value = 0.1768985792936573
$ws.Range("Q14").Value = 0.1289635301483333
$ws.Range("R14").Value = 0.77378118089
$ws.Range("S14").Value = 0.001500438288422752
$ws.Range("T14").Value = 0.001214079496086637

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1188766666666667
$ws.Range("H15").Value = 0.35663
$ws.Range("I15").Value = 0.006863138759702679
$ws.Range("J15").Value = 0.006863138759702677
$ws.Range("O15").Value = 0.2203773817607929
$ws.Range("P15").Value = 0.267477505734296
$ws.Range("Q15").Value = 0.1299985694722222
$ws.Range("R15").Value = 1.16998712525
$ws.Range("S15").Value = 0.001512480550524292
$ws.Range("T15").Value = 0.001835735236953642

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1188766666666667
$ws.Range("H16").Value = 0.35663
$ws.Range("I16").Value = 0.006863138759702679
$ws.Range("J16").Value = 0.006863138759702677
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.1406063333333333
$ws.Range("N16").Value = 0.421819
$ws.Range("O16").Value = 0.02833543913888328
$ws.Range("P16").Value = 0.03439142676166795
$ws.Range("Q16").Value = 0.01671481221888889
$ws.Range("R16").Value = 0.15043330997
$ws.Range("S16").Value = 0.0001944700506272661
$ws.Range("T16").Value = 0.0002360331340094793

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.1188766666666667
$ws.Range("H17").Value = 0.35663
$ws.Range("I17").Value = 0.006863138759702679
$ws.Range("J17").Value = 0.006863138759702677
$ws.Range("M17").Value = 1.5365345
$ws.Range("N17").Value = 3.073069
$ws.Range("O17").Value = 0.3096473592432615
$ws.Range("P17").Value = 0.2505511308097838
$ws.Range("Q17").Value = 0.1826580995783334
$ws.Range("R17").Value = 1.09594859747
$ws.Range("S17").Value = 0.002125152793062008
$ws.Range("T17").Value = 0.001719567177147963

# Row 18
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 0.6666666666666666
$ws.Range("G18").Value = 0.1188766666666667
$ws.Range("H18").Value = 0.35663
$ws.Range("I18").Value = 0.006863138759702679
$ws.Range("J18").Value = 0.006863138759702677
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = 0.6666666666666666
$ws.Range("M18").Value = 0.3224876666666667
$ws.Range("N18").Value = 0.967463
$ws.Range("O18").Value = 0.06498874862351253
$ws.Range("P18").Value = 0.07887845950306545
$ws.Range("Q18").Value = 0.03833625885444444
$ws.Range("R18").Value = 0.34502632969
$ws.Range("S18").Value = 0.000446026799622603
$ws.Range("T18").Value = 0.0005413538127211265

# Row 19
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 0.6666666666666666
$ws.Range("G19").Value = 0.1188766666666667
$ws.Range("H19").Value = 0.35663
$ws.Range("I19").Value = 0.006863138759702679
$ws.Range("J19").Value = 0.006863138759702677
$ws.Range("M19").Value = 0.7841693333333334
$ws.Range("N19").Value = 2.352508
$ws.Range("O19").Value = 0.1580283184440151
$ws.Range("P19").Value = 0.1918028978975294
$ws.Range("Q19").Value = 0.0932194364488889
$ws.Range("R19").Value = 0.8389749280400001
$ws.Range("S19").Value = 0.001084570277443758
